$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: unicode character class example
$ws.Range("B13").Value = "[\x{0800}-\x{FFFF}]"
$ws.Range("M13").Value = "unicode characters from 0800 to FFFF"

# Row 15: Chipseal example
$ws.Range("B15").Value = ".*(Chipseal){1}"
$ws.Range("D15").Value = "Renewal - Chipseal"
$ws.Range("M15").Value = "End with Chipseal"

$ws.Range("M16").Select()
